# Apply the TPM-recomputation update to the Fgf17-Fgfr4 LR-pairs sheet.
# - Adds a new "Resolving-Mac" sending-cluster block (rows 10-13)
# - Refreshes the recomputed numeric columns (H..T) on the existing rows (2-9)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Refresh recomputed numeric values on existing data rows 2-9 ---
# Row 2
$ws.Range("H2").Value = [double]"4.292218999999999"
$ws.Range("I2").Value = [double]"0.7311024264480693"
$ws.Range("J2").Value = [double]"0.7311024264480693"
$ws.Range("K2").Value = [double]"2"
$ws.Range("L2").Value = [double]"0.6666666666666666"
$ws.Range("M2").Value = [double]"0.3659943333333333"
$ws.Range("N2").Value = [double]"1.097983"
$ws.Range("O2").Value = [double]"0.006726051721149161"
$ws.Range("P2").Value = [double]"0.006726051721149162"
$ws.Range("Q2").Value = [double]"0.5236426104752221"
$ws.Range("R2").Value = [double]"4.712783494276999"
$ws.Range("S2").Value = [double]"0.004917432733747365"
$ws.Range("T2").Value = [double]"0.004917432733747365"

# Row 3
$ws.Range("H3").Value = [double]"4.292218999999999"
$ws.Range("I3").Value = [double]"0.7311024264480693"
$ws.Range("J3").Value = [double]"0.7311024264480693"
$ws.Range("O3").Value = [double]"0.001678071748088335"
$ws.Range("P3").Value = [double]"0.001678071748088335"
$ws.Range("S3").Value = [double]"0.001226842326781335"
$ws.Range("T3").Value = [double]"0.001226842326781335"

# Row 4
$ws.Range("H4").Value = [double]"4.292218999999999"
$ws.Range("I4").Value = [double]"0.7311024264480693"
$ws.Range("J4").Value = [double]"0.7311024264480693"
$ws.Range("M4").Value = [double]"53.897087"
$ws.Range("N4").Value = [double]"161.691261"
$ws.Range("O4").Value = [double]"0.9904923704135933"
$ws.Range("P4").Value = [double]"0.9904923704135934"
$ws.Range("Q4").Value = [double]"77.11270028868432"
$ws.Range("R4").Value = [double]"694.0143025981589"
$ws.Range("S4").Value = [double]"0.7241513753876779"
$ws.Range("T4").Value = [double]"0.724151375387678"

# Row 5
$ws.Range("H5").Value = [double]"4.292218999999999"
$ws.Range("I5").Value = [double]"0.7311024264480693"
$ws.Range("J5").Value = [double]"0.7311024264480693"
$ws.Range("K5").Value = [double]"2"
$ws.Range("L5").Value = [double]"0.6666666666666666"
$ws.Range("M5").Value = [double]"0.06004666666666667"
$ws.Range("N5").Value = [double]"0.18014"
$ws.Range("O5").Value = [double]"0.001103506117169219"
$ws.Range("P5").Value = [double]"0.001103506117169219"
$ws.Range("Q5").Value = [double]"0.08591114785111111"
$ws.Range("R5").Value = [double]"0.7732003306599999"
$ws.Range("S5").Value = [double]"0.0008067759998627032"
$ws.Range("T5").Value = [double]"0.0008067759998627033"

# Row 6
$ws.Range("G6").Value = [double]"0.4223926666666666"
$ws.Range("H6").Value = [double]"1.267178"
$ws.Range("I6").Value = [double]"0.2158410161600821"
$ws.Range("J6").Value = [double]"0.2158410161600822"
$ws.Range("K6").Value = [double]"2"
$ws.Range("L6").Value = [double]"0.6666666666666666"
$ws.Range("M6").Value = [double]"0.3659943333333333"
$ws.Range("N6").Value = [double]"1.097983"
$ws.Range("O6").Value = [double]"0.006726051721149161"
$ws.Range("P6").Value = [double]"0.006726051721149162"
$ws.Range("Q6").Value = [double]"0.1545933224415555"
$ws.Range("R6").Value = [double]"1.391339901974"
$ws.Range("S6").Value = [double]"0.001451757838238104"
$ws.Range("T6").Value = [double]"0.001451757838238105"

# Row 7
$ws.Range("G7").Value = [double]"0.4223926666666666"
$ws.Range("H7").Value = [double]"1.267178"
$ws.Range("I7").Value = [double]"0.2158410161600821"
$ws.Range("J7").Value = [double]"0.2158410161600822"
$ws.Range("O7").Value = [double]"0.001678071748088335"
$ws.Range("P7").Value = [double]"0.001678071748088335"
$ws.Range("Q7").Value = [double]"0.03856923758355556"
$ws.Range("R7").Value = [double]"0.347123138252"
$ws.Range("S7").Value = [double]"0.0003621967112969117"
$ws.Range("T7").Value = [double]"0.0003621967112969117"

# Row 8
$ws.Range("G8").Value = [double]"0.4223926666666666"
$ws.Range("H8").Value = [double]"1.267178"
$ws.Range("I8").Value = [double]"0.2158410161600821"
$ws.Range("J8").Value = [double]"0.2158410161600822"
$ws.Range("M8").Value = [double]"53.897087"
$ws.Range("N8").Value = [double]"161.691261"
$ws.Range("O8").Value = [double]"0.9904923704135933"
$ws.Range("P8").Value = [double]"0.9904923704135934"
$ws.Range("Q8").Value = [double]"22.76573430349533"
$ws.Range("R8").Value = [double]"204.891608731458"
$ws.Range("S8").Value = [double]"0.2137888797288784"
$ws.Range("T8").Value = [double]"0.2137888797288785"

# Row 9
$ws.Range("G9").Value = [double]"0.4223926666666666"
$ws.Range("H9").Value = [double]"1.267178"
$ws.Range("I9").Value = [double]"0.2158410161600821"
$ws.Range("J9").Value = [double]"0.2158410161600822"
$ws.Range("K9").Value = [double]"2"
$ws.Range("L9").Value = [double]"0.6666666666666666"
$ws.Range("M9").Value = [double]"0.06004666666666667"
$ws.Range("N9").Value = [double]"0.18014"
$ws.Range("O9").Value = [double]"0.001103506117169219"
$ws.Range("P9").Value = [double]"0.001103506117169219"
$ws.Range("Q9").Value = [double]"0.02536327165777778"
$ws.Range("R9").Value = [double]"0.22826944492"
$ws.Range("S9").Value = [double]"0.0002381818816686708"
$ws.Range("T9").Value = [double]"0.0002381818816686709"

# --- Append new rows 10-13: "Resolving-Mac" as a sending cluster ---
# Row 10
$ws.Range("A10").Value = "Resolving-Mac"
$ws.Range("B10").Value = "Fgf17"
$ws.Range("C10").Value = "Fgfr4"
$ws.Range("D10").Value = "ECs"
$ws.Range("E10").Value = [double]"1"
$ws.Range("F10").Value = [double]"0.3333333333333333"
$ws.Range("G10").Value = [double]"0.1038296666666667"
$ws.Range("H10").Value = [double]"0.311489"
$ws.Range("I10").Value = [double]"0.05305655739184852"
$ws.Range("J10").Value = [double]"0.05305655739184854"
$ws.Range("K10").Value = [double]"2"
$ws.Range("L10").Value = [double]"0.6666666666666666"
$ws.Range("M10").Value = [double]"0.3659943333333333"
$ws.Range("N10").Value = [double]"1.097983"
$ws.Range("O10").Value = [double]"0.006726051721149161"
$ws.Range("P10").Value = [double]"0.006726051721149162"
$ws.Range("Q10").Value = [double]"0.03800106963188889"
$ws.Range("R10").Value = [double]"0.342009626687"
$ws.Range("S10").Value = [double]"0.000356861149163692"
$ws.Range("T10").Value = [double]"0.0003568611491636922"

# Row 11
$ws.Range("A11").Value = "Resolving-Mac"
$ws.Range("B11").Value = "Fgf17"
$ws.Range("C11").Value = "Fgfr4"
$ws.Range("D11").Value = "FAPs"
$ws.Range("E11").Value = [double]"1"
$ws.Range("F11").Value = [double]"0.3333333333333333"
$ws.Range("G11").Value = [double]"0.1038296666666667"
$ws.Range("H11").Value = [double]"0.311489"
$ws.Range("I11").Value = [double]"0.05305655739184852"
$ws.Range("J11").Value = [double]"0.05305655739184854"
$ws.Range("K11").Value = [double]"1"
$ws.Range("L11").Value = [double]"0.3333333333333333"
$ws.Range("M11").Value = [double]"0.09131133333333334"
$ws.Range("N11").Value = [double]"0.273934"
$ws.Range("O11").Value = [double]"0.001678071748088335"
$ws.Range("P11").Value = [double]"0.001678071748088335"
$ws.Range("Q11").Value = [double]"0.009480825302888889"
$ws.Range("R11").Value = [double]"0.08532742772600001"
$ws.Range("S11").Value = [double]"8.903271001008834E-05"
$ws.Range("T11").Value = [double]"8.903271001008837E-05"

# Row 12
$ws.Range("A12").Value = "Resolving-Mac"
$ws.Range("B12").Value = "Fgf17"
$ws.Range("C12").Value = "Fgfr4"
$ws.Range("D12").Value = "MuSCs"
$ws.Range("E12").Value = [double]"1"
$ws.Range("F12").Value = [double]"0.3333333333333333"
$ws.Range("G12").Value = [double]"0.1038296666666667"
$ws.Range("H12").Value = [double]"0.311489"
$ws.Range("I12").Value = [double]"0.05305655739184852"
$ws.Range("J12").Value = [double]"0.05305655739184854"
$ws.Range("K12").Value = [double]"3"
$ws.Range("L12").Value = [double]"1"
$ws.Range("M12").Value = [double]"53.897087"
$ws.Range("N12").Value = [double]"161.691261"
$ws.Range("O12").Value = [double]"0.9904923704135933"
$ws.Range("P12").Value = [double]"0.9904923704135934"
$ws.Range("Q12").Value = [double]"5.596116577514334"
$ws.Range("R12").Value = [double]"50.365049197629"
$ws.Range("S12").Value = [double]"0.0525521152970369"
$ws.Range("T12").Value = [double]"0.05255211529703692"

# Row 13
$ws.Range("A13").Value = "Resolving-Mac"
$ws.Range("B13").Value = "Fgf17"
$ws.Range("C13").Value = "Fgfr4"
$ws.Range("D13").Value = "Resolving-Mac"
$ws.Range("E13").Value = [double]"1"
$ws.Range("F13").Value = [double]"0.3333333333333333"
$ws.Range("G13").Value = [double]"0.1038296666666667"
$ws.Range("H13").Value = [double]"0.311489"
$ws.Range("I13").Value = [double]"0.05305655739184852"
$ws.Range("J13").Value = [double]"0.05305655739184854"
$ws.Range("K13").Value = [double]"2"
$ws.Range("L13").Value = [double]"0.6666666666666666"
$ws.Range("M13").Value = [double]"0.06004666666666667"
$ws.Range("N13").Value = [double]"0.18014"
$ws.Range("O13").Value = [double]"0.001103506117169219"
$ws.Range("P13").Value = [double]"0.001103506117169219"
$ws.Range("Q13").Value = [double]"0.006234625384444445"
$ws.Range("R13").Value = [double]"0.05611162846000001"
$ws.Range("S13").Value = [double]"5.854823563784457E-05"
$ws.Range("T13").Value = [double]"5.854823563784459E-05"
